$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 4 datetime values (Correspond Handoff / Handback Datetime)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-19 06:34:08"
$wsZh.Range("H4").Value = "2016-03-19 06:34:28"

# de-de sheet: row 4 datetime values (Correspond Handoff / Handback Datetime)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-19 06:34:10"
$wsDe.Range("H4").Value = "2016-03-19 06:34:32"
